# Weekly update: insert a new daily price record for Mango
# (Vega Central Mapocho de Santiago) above the existing row 355,
# shifting all subsequent rows (old 355-417) down by one (new 356-418).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 355; this shifts rows 355:417 -> 356:418
$ws.Rows("355").Insert()

# Populate the newly inserted row 355 with the new record's data.
$ws.Range("A355").Value = 9
$ws.Range("B355").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C355").Value = "Metropolitana"
$ws.Range("D355").Value = 44694
$ws.Range("E355").Value = 13
$ws.Range("F355").Value = "Fruta"
$ws.Range("G355").Value = 100108
$ws.Range("H355").Value = "Tropicales y subtropicales"
$ws.Range("I355").Value = 100108002
$ws.Range("J355").Value = "Mango"
$ws.Range("K355").Value = "Sin especificar"
$ws.Range("L355").Value = "Primera"
$ws.Range("M355").Value = 480
$ws.Range("N355").Value = 7000
$ws.Range("O355").Value = 7500
$ws.Range("P355").Value = 7292
$ws.Range("Q355").Value = "$/bandeja 4 kilos"
$ws.Range("R355").Value = "Brasil"
$ws.Range("S355").Value = 1823
$ws.Range("T355").Value = 4
